$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
  2  = "  -0.53%  "
  3  = "  -0.90%  "
  4  = "  -0.03%  "
  5  = "  -0.48%  "
  6  = "  +0.40%  "
  8  = "  +0.72%  "
  9  = "  +5.06%  "
  10 = "  -1.49%  "
  11 = "  +2.14%  "
  12 = "  -0.95%  "
  13 = "  +0.14%  "
  14 = "  -5.48%  "
  15 = "  -0.49%  "
  16 = "  -0.30%  "
  18 = "  -0.39%  "
  19 = "  +1.42%  "
  20 = "  -3.21%  "
  21 = "  -1.04%  "
  22 = "  -0.38%  "
  23 = "  +0.19%  "
  24 = "  +0.02%  "
  25 = "  -1.01%  "
  26 = "  -0.84%  "
  27 = "  +1.28%  "
  28 = "  -0.68%  "
  29 = "  -0.26%  "
  30 = "  -2.64%  "
  31 = "  -1.78%  "
  32 = "  -0.01%  "
  33 = "  -4.73%  "
  34 = "  -2.55%  "
  35 = "  +1.48%  "
  36 = "  +3.81%  "
  37 = "  +1.62%  "
  38 = "  -2.60%  "
  39 = "  -1.20%  "
  40 = "  -1.04%  "
  41 = "  +0.20%  "
  42 = "  -2.28%  "
  43 = "  +0.77%  "
  44 = "  -3.10%  "
  45 = "  -1.61%  "
  46 = "  -1.75%  "
  47 = "  -1.13%  "
  48 = "  +1.45%  "
  49 = "  -0.93%  "
  50 = "  -4.35%  "
  51 = "  +0.04%  "
}

foreach ($row in $updates.Keys) {
  $ws.Range("E$row").Value = $updates[$row]
}
